$d = $word.ActiveDocument

# Change the Team ID: "TMID22194" -> "TMID11232", and remove one leading space
# before "NM2023" in that cell.
$d.Content.Find.Execute("TMID22194", $true, $false, $false, $false, $false,
                         $true, 1, $false, "TMID11232", 2) | Out-Null

$d.Content.Find.Execute("          NM2023", $true, $false, $false, $false, $false,
                         $true, 1, $false, "         NM2023", 2) | Out-Null
